$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update closing price (종가) and 5-day return (5일수익률) for row 2 (TSM)
$ws.Range("D2").Value = 291.51
$ws.Range("F2").Value = 5.05

# Update MACRO_SCORE column (N) for rows 2 through 6
$ws.Range("N2").Value = 85.83574689470727
$ws.Range("N3").Value = 85.83574689470727
$ws.Range("N4").Value = 85.83574689470727
$ws.Range("N5").Value = 85.83574689470727
$ws.Range("N6").Value = 85.83574689470727
